# Update covid country stats & re-sorted "provincias Spain" rows (paises.xlsx)
# - refreshes the "datos actualizados" timestamp (19:05 -> 20:05)
# - refreshes Casos totales/Nuevos casos/Casos activos/Recuperados/Casos criticos/Muertes hoy/Muertes
#   for the countries whose counters moved between this export and the previous one
# - some countries changed total-case rank and therefore swapped sorted rows; for those rows
#   both the country name (col A) and its stats needed to move down/up to the new row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp cell (A1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 31 de Mayo de 2020 a las 20:05"  # col A

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 1826394  # col B
$ws.Cells.Item(4, 3).Value = 9574  # col C
$ws.Cells.Item(4, 4).Value = 538587  # col D
$ws.Cells.Item(4, 5).Value = 1181926  # col E
$ws.Cells.Item(4, 7).Value = 324  # col G
$ws.Cells.Item(4, 8).Value = 105881  # col H

# Row 7: España
$ws.Cells.Item(7, 2).Value = 286509  # col B
$ws.Cells.Item(7, 3).Value = 201  # col C
$ws.Cells.Item(7, 5).Value = 62424  # col E
$ws.Cells.Item(7, 7).Value = 2  # col G
$ws.Cells.Item(7, 8).Value = 27127  # col H

# Row 12: Alemania
$ws.Cells.Item(12, 2).Value = 183442  # col B
$ws.Cells.Item(12, 3).Value = 148  # col C
$ws.Cells.Item(12, 5).Value = 9640  # col E

# Row 13: Turquia
$ws.Cells.Item(13, 2).Value = 163942  # col B
$ws.Cells.Item(13, 3).Value = 839  # col C
$ws.Cells.Item(13, 4).Value = 127973  # col D
$ws.Cells.Item(13, 5).Value = 31429  # col E
$ws.Cells.Item(13, 7).Value = 25  # col G
$ws.Cells.Item(13, 8).Value = 4540  # col H

# Row 16: Chile
$ws.Cells.Item(16, 2).Value = 99688  # col B
$ws.Cells.Item(16, 3).Value = 4830  # col C
$ws.Cells.Item(16, 4).Value = 42727  # col D
$ws.Cells.Item(16, 5).Value = 55907  # col E
$ws.Cells.Item(16, 7).Value = 57  # col G
$ws.Cells.Item(16, 8).Value = 1054  # col H

# Row 44: Israel
$ws.Cells.Item(44, 2).Value = 17071  # col B
$ws.Cells.Item(44, 3).Value = 59  # col C
$ws.Cells.Item(44, 5).Value = 1974  # col E
$ws.Cells.Item(44, 7).Value = 1  # col G
$ws.Cells.Item(44, 8).Value = 285  # col H

# Row 65: Marruecos
$ws.Cells.Item(65, 2).Value = 7807  # col B
$ws.Cells.Item(65, 3).Value = 27  # col C
$ws.Cells.Item(65, 4).Value = 5459  # col D
$ws.Cells.Item(65, 5).Value = 2143  # col E
$ws.Cells.Item(65, 7).Value = 1  # col G
$ws.Cells.Item(65, 8).Value = 205  # col H

# Row 101: Sri Lanka
$ws.Cells.Item(101, 2).Value = 1633  # col B
$ws.Cells.Item(101, 3).Value = 20  # col C
$ws.Cells.Item(101, 5).Value = 822  # col E

# Row 120: Paraguay
$ws.Cells.Item(120, 2).Value = 986  # col B
$ws.Cells.Item(120, 3).Value = 22  # col C
$ws.Cells.Item(120, 4).Value = 477  # col D
$ws.Cells.Item(120, 5).Value = 498  # col E

# Row 127: Madagascar -> Republica del Chad
$ws.Cells.Item(127, 1).Value = "Republica del Chad"  # col A
$ws.Cells.Item(127, 2).Value = 778  # col B
$ws.Cells.Item(127, 3).Value = 19  # col C
$ws.Cells.Item(127, 4).Value = 491  # col D
$ws.Cells.Item(127, 5).Value = 222  # col E
$ws.Cells.Item(127, 8).Value = 65  # col H

# Row 128: Principado de Andorra -> Madagascar
$ws.Cells.Item(128, 1).Value = "Madagascar"  # col A
$ws.Cells.Item(128, 2).Value = 771  # col B
$ws.Cells.Item(128, 3).Value = 13  # col C
$ws.Cells.Item(128, 4).Value = 168  # col D
$ws.Cells.Item(128, 5).Value = 597  # col E
$ws.Cells.Item(128, 8).Value = 6  # col H

# Row 129: Nicaragua -> Principado de Andorra
$ws.Cells.Item(129, 1).Value = "Principado de Andorra"  # col A
$ws.Cells.Item(129, 2).Value = 764  # col B
$ws.Cells.Item(129, 4).Value = 692  # col D
$ws.Cells.Item(129, 5).Value = 21  # col E
$ws.Cells.Item(129, 8).Value = 51  # col H

# Row 130: Republica del Chad -> Nicaragua
$ws.Cells.Item(130, 1).Value = "Nicaragua"  # col A
$ws.Cells.Item(130, 4).Value = 370  # col D
$ws.Cells.Item(130, 5).Value = 354  # col E
$ws.Cells.Item(130, 8).Value = 35  # col H

# Row 135: Jamaica -> Congo
$ws.Cells.Item(135, 1).Value = "Congo"  # col A
$ws.Cells.Item(135, 2).Value = 611  # col B
$ws.Cells.Item(135, 3).Value = 40  # col C
$ws.Cells.Item(135, 4).Value = 179  # col D
$ws.Cells.Item(135, 5).Value = 412  # col E
$ws.Cells.Item(135, 7).Value = 1  # col G
$ws.Cells.Item(135, 8).Value = 20  # col H

# Row 136: Congo -> Jamaica
$ws.Cells.Item(136, 1).Value = "Jamaica"  # col A
$ws.Cells.Item(136, 2).Value = 581  # col B
$ws.Cells.Item(136, 3).Value = 6  # col C
$ws.Cells.Item(136, 4).Value = 290  # col D
$ws.Cells.Item(136, 5).Value = 282  # col E
$ws.Cells.Item(136, 8).Value = 9  # col H

# Row 144: Togo -> Cabo Verde
$ws.Cells.Item(144, 1).Value = "Cabo Verde"  # col A
$ws.Cells.Item(144, 2).Value = 435  # col B
$ws.Cells.Item(144, 3).Value = 14  # col C
$ws.Cells.Item(144, 4).Value = 193  # col D
$ws.Cells.Item(144, 5).Value = 238  # col E
$ws.Cells.Item(144, 8).Value = 4  # col H

# Row 145: Cabo Verde -> Togo
$ws.Cells.Item(145, 1).Value = "Togo"  # col A
$ws.Cells.Item(145, 2).Value = 433  # col B
$ws.Cells.Item(145, 4).Value = 206  # col D
$ws.Cells.Item(145, 5).Value = 214  # col E
$ws.Cells.Item(145, 8).Value = 13  # col H

# Row 146: Uganda
$ws.Cells.Item(146, 2).Value = 417  # col B
$ws.Cells.Item(146, 3).Value = 4  # col C
$ws.Cells.Item(146, 5).Value = 345  # col E

# Row 154: Suazilandia
$ws.Cells.Item(154, 2).Value = 285  # col B
$ws.Cells.Item(154, 3).Value = 2  # col C
$ws.Cells.Item(154, 4).Value = 189  # col D
$ws.Cells.Item(154, 5).Value = 94  # col E

# Row 158: Birmania
$ws.Cells.Item(158, 4).Value = 138  # col D
$ws.Cells.Item(158, 5).Value = 80  # col E

# Row 200: Santa Lucia -> Belice
$ws.Cells.Item(200, 1).Value = "Belice"  # col A
$ws.Cells.Item(200, 4).Value = 16  # col D
$ws.Cells.Item(200, 8).Value = 2  # col H

# Row 201: Belice -> Santa Lucia
$ws.Cells.Item(201, 1).Value = "Santa Lucia"  # col A
$ws.Cells.Item(201, 4).Value = 18  # col D
$ws.Cells.Item(201, 8).Value = 0  # col H

# Row 213: Papua Nueva Guinea -> Islas Virgenes Britanicas
$ws.Cells.Item(213, 1).Value = "Islas Virgenes Britanicas"  # col A
$ws.Cells.Item(213, 4).Value = 7  # col D
$ws.Cells.Item(213, 8).Value = 1  # col H

# Row 214: Islas Virgenes Britanicas -> Papua Nueva Guinea
$ws.Cells.Item(214, 1).Value = "Papua Nueva Guinea"  # col A
$ws.Cells.Item(214, 4).Value = 8  # col D
$ws.Cells.Item(214, 8).Value = 0  # col H
